$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table and Chart")

# Update "Split Count" column (C) value for row 4 (13 -> 15)
$ws.Range("C4").Value = 15

# Update "Split Sort Conut" column (B) values
$ws.Range("B5").Value = 5754
$ws.Range("B6").Value = 11914
$ws.Range("B7").Value = 25130
$ws.Range("B8").Value = 28766

# Update "Split Count" column (C) values
$ws.Range("C5").Value = 2226
$ws.Range("C6").Value = 4435
$ws.Range("C7").Value = 8084
$ws.Range("C8").Value = 9900
